$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Marking" row values (row 11)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Fix the "Total" row values (row 12)
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "74 / 112"
